$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update count column (C) values
$ws.Range("C2").Value = 8
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 10
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 13
$ws.Range("C9").Value = 6
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 6
$ws.Range("C12").Value = 8
$ws.Range("C13").Value = 7
$ws.Range("C14").Value = 8
$ws.Range("C15").Value = 6
$ws.Range("C16").Value = 1
$ws.Range("C17").Value = 8

# Update recognized word text column (B)
$ws.Range("B8").Value = "<november>"
$ws.Range("B17").Value = "<wouk>"
